$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = "aluminium"
$ws.Range("B19").Value = 0.002

$ws.Range("A20").Value = "Polyéthylène"
$ws.Range("B20").Value = 0.002

$ws.Range("A21").Select()
